# "For PCBWAY March 03"
#
# The BOM's "Line #" column (A) was blank for every component row; it is
# now filled in with an explicit empty text value (equivalent to typing a
# lone leading apostrophe into the cell, i.e. an empty "quote-prefixed"
# text entry) for every data row (2-12).
#
# Two "Manufacturer Lifecycle 1" cells (J4, J11) that said "Unknown" are
# cleared to the same kind of empty text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lineNumberCells = @("A2", "A3", "A4", "A5", "A6", "A7", "A8", "A9", "A10", "A11", "A12")
foreach ($cellRef in $lineNumberCells) {
    $ws.Range($cellRef).Value = "'"
}

$unknownCells = @("J4", "J11")
foreach ($cellRef in $unknownCells) {
    $ws.Range($cellRef).Value = "'"
}
